$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Date" column contains a mis-formatted date string ("6-12-2011-12")
# for every data row. The NBA stats feed was off by a day when the data was
# scraped, so the value should actually read "2012-06-12".
#
# Locate the "Date" header and the extent of the data dynamically rather
# than hard-coding the column/row numbers.
$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $used.Row + $used.Rows.Count - 1
$firstCol = $used.Column
$lastCol = $used.Column + $used.Columns.Count - 1

$dateCol = 0
for ($c = $firstCol; $c -le $lastCol; $c++) {
    if ($ws.Cells.Item($firstRow, $c).Text -eq "Date") {
        $dateCol = $c
    }
}

if ($dateCol -gt 0) {
    for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, $dateCol)
        if ($cell.Text -eq "6-12-2011-12") {
            # Assigning the literal string "2012-06-12" directly would make
            # Excel auto-detect it as a date and silently convert the cell
            # to a date-formatted serial number, which is not the desired
            # result - the corrected value must remain plain text, just
            # like the original. Building it via a formula that evaluates
            # to text, then collapsing that formula to its static result in
            # place (copy / paste-special values only), keeps it as plain
            # text without introducing a new number format or cell style.
            $cell.Formula = '="2012-06-12"'
            $cell.Copy()
            $cell.PasteSpecial(-4163)
        }
    }
}

$excel.CutCopyMode = 0
